$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "DLM_SYSTEM_STATUS"
$ws.Range("B14").Value = "PROD"
$ws.Range("C14").Value = "Database Logging Module data system status - this configuration option determines which messages are logged in the database"
$ws.Range("D14").Formula = '="INSERT INTO CC_CONFIG_OPTIONS ("&$A$1&", "&$B$1&", "&$C$1&") VALUES (''"&SUBSTITUTE(A14, "''", "''''")&"'', ''"&SUBSTITUTE(B14, "''", "''''")&"'', ''"&SUBSTITUTE(C14, "''", "''''")&"'');"'

$ws.Range("D14").Select()
